# "Development Center in Extension" - add a DEV-MENU branch with
# develop.g.* graph-tool entries, and re-home the old "图管理"
# (graph management) entry under the new parent as a legacy item.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Build out the style palette we'll need for the new rows by
#    copying (format-only) from existing cells that already carry
#    the exact xf we want, then writing the cell content on top.
# ---------------------------------------------------------------

# Style 8 (id columns A/B) -> rows 6-11
$ws.Range("A5:B5").Copy()
$ws.Range("A6:B11").PasteSpecial(-4122)

# Style 12 (table-name column C) -> rows 6-11
$ws.Range("C5").Copy()
$ws.Range("C6:C11").PasteSpecial(-4122)

# Style 9 (order column D, center) -> rows 7-11 (row 6 gets its own below)
$ws.Range("D5").Copy()
$ws.Range("D7:D11").PasteSpecial(-4122)

# Style 16 (level column E, center) -> rows 7-11 (row 6 gets its own below)
$ws.Range("E5").Copy()
$ws.Range("E7:E11").PasteSpecial(-4122)

# Style 14 (name column F) -> rows 6-11
$ws.Range("F5").Copy()
$ws.Range("F6:F11").PasteSpecial(-4122)

# Style 13 (text/icon columns G/H) -> rows 6-11
$ws.Range("G5:H5").Copy()
$ws.Range("G6:H11").PasteSpecial(-4122)

# Style 15 (uri column I) -> rows 7-11 (row 6 gets its own below)
$ws.Range("I5").Copy()
$ws.Range("I7:I11").PasteSpecial(-4122)

# Row 6 is the new parent "develop.g" row - it reuses the header-ish
# highlighted xfs already present on row 4 (fontId3/fillId2) for its
# order/level cells.
$ws.Range("D4").Copy()
$ws.Range("D6").PasteSpecial(-4122)

$ws.Range("D5").Copy()
$ws.Range("E6").PasteSpecial(-4122)

# Row 6's uri cell ("EXPAND") needs a brand-new style: same box as the
# other uri cells (left aligned, bordered, no fill) but in bold red
# 16pt text. Clone the left-aligned bordered xf, then recolor the font;
# the engine will mint the missing font+xf automatically.
$ws.Range("G5").Copy()
$ws.Range("I6").PasteSpecial(-4122)
$ws.Range("I6").Font.Size = 16
$ws.Range("I6").Font.Color = 255

# ---------------------------------------------------------------
# 2. Re-point the existing row 5 at the new parent chain. The old
#    "图管理"/SIDE-MENU record becomes a legacy DEV-MENU leaf.
# ---------------------------------------------------------------

$ws.Range("A5").Value = "7861a7ae-2cb0-49e6-9a57-01d0ecf6ebed"
$ws.Range("C5").Value = "DEV-MENU"
$ws.Range("D5").Value = 2000
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = "zero.develop.graphic"
$ws.Range("G5").Value = "「旧」图管理"
$ws.Range("H5").Value = "share-alt"
$ws.Range("I5").Value = "/epic/topology"
# B5's key now mirrors the new parent's id via formula
$ws.Range("B5").Formula = "=A$6"

# ---------------------------------------------------------------
# 3. Row 6 - the new "develop.g" parent menu entry (no B/key, no J).
# ---------------------------------------------------------------

$ws.Range("A6").Value = "a52f5024-a481-4777-8b04-a8e27a5af72a"
$ws.Range("C6").Value = "DEV-MENU"
$ws.Range("D6").Value = 50000
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = "develop.g"
$ws.Range("G6").Value = "图管理"
$ws.Range("H6").Value = "reconciliation"
$ws.Range("I6").Value = "EXPAND"

# ---------------------------------------------------------------
# 4. Rows 7-11 - develop.g.* leaf entries, all keyed off A$6 via a
#    shared formula in column B.
# ---------------------------------------------------------------

$ws.Range("A7").Value = "70757329-ff2f-4647-a710-f11b3dea7975"
$ws.Range("C7").Value = "DEV-MENU"
$ws.Range("D7").Value = 1005
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = "develop.g.relation"
$ws.Range("G7").Value = "关系图"
$ws.Range("H7").Value = "block"
$ws.Range("I7").Value = "/g/relation"

$ws.Range("A8").Value = "dc36e977-e6a6-4f21-9219-cee6556807b6"
$ws.Range("C8").Value = "DEV-MENU"
$ws.Range("D8").Value = 1010
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = "develop.g.tree"
$ws.Range("G8").Value = "树型图"
$ws.Range("H8").Value = "cluster"
$ws.Range("I8").Value = "/g/tree"

$ws.Range("A9").Value = "48873c15-49fa-41f3-9e7f-844bcc1256f7"
$ws.Range("C9").Value = "DEV-MENU"
$ws.Range("D9").Value = 1015
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = "develop.g.workflow"
$ws.Range("G9").Value = "流程图"
$ws.Range("H9").Value = "gateway"
$ws.Range("I9").Value = "/g/workflow"

$ws.Range("A10").Value = "4cbd5021-2a70-4dab-bdef-91ca74d53fa6"
$ws.Range("C10").Value = "DEV-MENU"
$ws.Range("D10").Value = 1020
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = "develop.g.management"
$ws.Range("G10").Value = "图管理/发布"
$ws.Range("H10").Value = "thunderbolt"
$ws.Range("I10").Value = "/g/management"

$ws.Range("A11").Value = "b189936e-ac95-4cb6-803a-f6402ded6caa"
$ws.Range("C11").Value = "DEV-MENU"
$ws.Range("D11").Value = 1025
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = "develop.g.designer"
$ws.Range("G11").Value = "图设计器"
$ws.Range("H11").Value = "edit"
$ws.Range("I11").Value = "/g/designer"

# B7 mirrors B5's formula as its own standalone formula, then B8:B11
# share a single formula group (matches the author's edit pattern).
$ws.Range("B7").Formula = "=A$6"
$ws.Range("B8:B11").Formula = "=A$6"

# ---------------------------------------------------------------
# 5. Selection, matching the author's last recorded cursor position.
# ---------------------------------------------------------------
$ws.Range("F7").Select()
